$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-24: id (col B), speaker_variant (col C)
# Column D (is_prefered) is cleared for every row in this range.
$rows = @(
    @{ Row = 2;  Id = "#polinos";       Name = "Polinos" },
    @{ Row = 3;  Id = "#montaan";       Name = "Montaan" },
    @{ Row = 4;  Id = "#fillida";       Name = "Fillida" },
    @{ Row = 5;  Id = "#soliman";       Name = "Soliman" },
    @{ Row = 6;  Id = "#armida";        Name = "Armida" },
    @{ Row = 7;  Id = "#acastus,-echo"; Name = "Acastus, Echo" },
    @{ Row = 8;  Id = "#cleandra";      Name = "Cleandra" },
    @{ Row = 9;  Id = "#ierahim";       Name = "Ierahim" },
    @{ Row = 10; Id = "#rein-out";      Name = "Rein out" },
    @{ Row = 11; Id = "#dares";         Name = "Dares" },
    @{ Row = 12; Id = "#cupido";        Name = "Cupido" },
    @{ Row = 13; Id = "#mustaffa";      Name = "Mustaffa" },
    @{ Row = 14; Id = "#tisbe";         Name = "Tisbe" },
    @{ Row = 15; Id = "#hydraöt";       Name = "Hydraöt" },
    @{ Row = 16; Id = "#zim";           Name = "Zim" },
    @{ Row = 17; Id = "#silvaan";       Name = "Silvaan" },
    @{ Row = 18; Id = "#reinout";       Name = "Reinout" },
    @{ Row = 19; Id = "#aurora";        Name = "Aurora" },
    @{ Row = 20; Id = "#geeandra";      Name = "Geeandra" },
    @{ Row = 21; Id = "#geweld";        Name = "Geweld" },
    @{ Row = 22; Id = "#acastus";       Name = "Acastus" },
    @{ Row = 23; Id = "#darfs";         Name = "Darfs" },
    @{ Row = 24; Id = "#ibrahim";       Name = "Ibrahim" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Id
    $ws.Cells.Item($r.Row, 3).Value = $r.Name
    $ws.Cells.Item($r.Row, 4).Value = ""
}
